$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 21: date 44438 -> 44463, quality Primera -> Especial, volume 100 -> 150
# (N21/O21/P21/S21 remain unchanged)
# ------------------------------------------------------------------
$ws.Range("D21").Value = 44463
$ws.Range("L21").Value = "Especial"
$ws.Range("M21").Value = 150

# ------------------------------------------------------------------
# Row 22: date 44461 -> 44463, quality Especial -> Primera, volume 150 -> 100,
# prices 30000 -> 26000, price/kg 3000 -> 2600
# ------------------------------------------------------------------
$ws.Range("D22").Value = 44463
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 26000
$ws.Range("O22").Value = 26000
$ws.Range("P22").Value = 26000
$ws.Range("S22").Value = 2600

# ------------------------------------------------------------------
# Row 23: date 44461 -> 44438, prices 25000 -> 30000, price/kg 2500 -> 3000
# (quality stays Primera, volume stays 100)
# ------------------------------------------------------------------
$ws.Range("D23").Value = 44438
$ws.Range("N23").Value = 30000
$ws.Range("O23").Value = 30000
$ws.Range("P23").Value = 30000
$ws.Range("S23").Value = 3000

# ------------------------------------------------------------------
# New row 24 (copy of former row22 record: Especial, 44461, 30000)
# ------------------------------------------------------------------
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44461
$ws.Range("D24").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100107
$ws.Range("H24").Value = "Otros"
$ws.Range("I24").Value = 100107002
$ws.Range("J24").Value = "Chirimoya"
$ws.Range("K24").Value = "Cultivar IV Región"
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 150
$ws.Range("N24").Value = 30000
$ws.Range("O24").Value = 30000
$ws.Range("P24").Value = 30000
$ws.Range("Q24").Value = "$/bandeja 10 kilos"
$ws.Range("R24").Value = "Provincia de Limarí"
$ws.Range("S24").Value = 3000
$ws.Range("T24").Value = 10

# ------------------------------------------------------------------
# New row 25 (copy of former row23 record: Primera, 44461, 25000)
# ------------------------------------------------------------------
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44461
$ws.Range("D25").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100107
$ws.Range("H25").Value = "Otros"
$ws.Range("I25").Value = 100107002
$ws.Range("J25").Value = "Chirimoya"
$ws.Range("K25").Value = "Cultivar IV Región"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 25000
$ws.Range("O25").Value = 25000
$ws.Range("P25").Value = 25000
$ws.Range("Q25").Value = "$/bandeja 10 kilos"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 2500
$ws.Range("T25").Value = 10

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
